# Adds rows 123-126 (IPL matches played on 2023-05-16 and 2023-05-17)
# to the batting/bowling team-innings stats sheet, per the "Update 21 May 2023" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 123
$ws.Range("A123").Value = "'2023-05-16"
$ws.Range("B123").Value = "Lucknow"
$ws.Range("C123").Value = "1st innings"
$ws.Range("D123").Value = "LSG"
$ws.Range("E123").Value = "MI"
$ws.Range("F123").Value = 1
$ws.Range("G123").Value = 167
$ws.Range("H123").Value = 34
$ws.Range("I123").Value = 77
$ws.Range("J123").Value = 56
$ws.Range("K123").Value = 4
$ws.Range("L123").Value = 2
$ws.Range("M123").Value = 2
$ws.Range("N123").Value = 0
$ws.Range("O123").Value = 125
$ws.Range("P123").Value = 36
$ws.Range("Q123").Value = 63
$ws.Range("R123").Value = 26
$ws.Range("S123").Value = "2nd innings"
$ws.Range("T123").Value = "LSG"
$ws.Range("U123").Value = 157
$ws.Range("V123").Value = 55
$ws.Range("W123").Value = 71
$ws.Range("X123").Value = 31
$ws.Range("Y123").Value = 5
$ws.Range("Z123").Value = 0
$ws.Range("AA123").Value = 3
$ws.Range("AB123").Value = 2
$ws.Range("AC123").Value = 20
$ws.Range("AD123").Value = 6
$ws.Range("AE123").Value = 10
$ws.Range("AF123").Value = 4

# Row 124
$ws.Range("A124").Value = "'2023-05-16"
$ws.Range("B124").Value = "Lucknow"
$ws.Range("C124").Value = "2nd innings"
$ws.Range("D124").Value = "MI"
$ws.Range("E124").Value = "LSG"
$ws.Range("F124").Value = 1
$ws.Range("G124").Value = 157
$ws.Range("H124").Value = 55
$ws.Range("I124").Value = 71
$ws.Range("J124").Value = 31
$ws.Range("K124").Value = 5
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 3
$ws.Range("N124").Value = 2
$ws.Range("O124").Value = 130
$ws.Range("P124").Value = 39
$ws.Range("Q124").Value = 62
$ws.Range("R124").Value = 29
$ws.Range("S124").Value = "1st innings"
$ws.Range("T124").Value = "MI"
$ws.Range("U124").Value = 167
$ws.Range("V124").Value = 34
$ws.Range("W124").Value = 77
$ws.Range("X124").Value = 56
$ws.Range("Y124").Value = 4
$ws.Range("Z124").Value = 2
$ws.Range("AA124").Value = 2
$ws.Range("AB124").Value = 0
$ws.Range("AC124").Value = 20
$ws.Range("AD124").Value = 6
$ws.Range("AE124").Value = 10
$ws.Range("AF124").Value = 4

# Row 125
$ws.Range("A125").Value = "'2023-05-17"
$ws.Range("B125").Value = "Dharamsala"
$ws.Range("C125").Value = "1st innings"
$ws.Range("D125").Value = "DC"
$ws.Range("E125").Value = "PBKS"
$ws.Range("F125").Value = 1
$ws.Range("G125").Value = 208
$ws.Range("H125").Value = 60
$ws.Range("I125").Value = 91
$ws.Range("J125").Value = 57
$ws.Range("K125").Value = 2
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = 2
$ws.Range("N125").Value = 0
$ws.Range("O125").Value = 125
$ws.Range("P125").Value = 37
$ws.Range("Q125").Value = 62
$ws.Range("R125").Value = 26
$ws.Range("S125").Value = "2nd innings"
$ws.Range("T125").Value = "DC"
$ws.Range("U125").Value = 188
$ws.Range("V125").Value = 44
$ws.Range("W125").Value = 87
$ws.Range("X125").Value = 57
$ws.Range("Y125").Value = 8
$ws.Range("Z125").Value = 1
$ws.Range("AA125").Value = 3
$ws.Range("AB125").Value = 4
$ws.Range("AC125").Value = 20
$ws.Range("AD125").Value = 6
$ws.Range("AE125").Value = 10
$ws.Range("AF125").Value = 4

# Row 126
$ws.Range("A126").Value = "'2023-05-17"
$ws.Range("B126").Value = "Dharamsala"
$ws.Range("C126").Value = "2nd innings"
$ws.Range("D126").Value = "PBKS"
$ws.Range("E126").Value = "DC"
$ws.Range("F126").Value = 1
$ws.Range("G126").Value = 188
$ws.Range("H126").Value = 44
$ws.Range("I126").Value = 87
$ws.Range("J126").Value = 57
$ws.Range("K126").Value = 8
$ws.Range("L126").Value = 1
$ws.Range("M126").Value = 3
$ws.Range("N126").Value = 4
$ws.Range("O126").Value = 124
$ws.Range("P126").Value = 36
$ws.Range("Q126").Value = 61
$ws.Range("R126").Value = 27
$ws.Range("S126").Value = "1st innings"
$ws.Range("T126").Value = "PBKS"
$ws.Range("U126").Value = 208
$ws.Range("V126").Value = 60
$ws.Range("W126").Value = 91
$ws.Range("X126").Value = 57
$ws.Range("Y126").Value = 2
$ws.Range("Z126").Value = 0
$ws.Range("AA126").Value = 2
$ws.Range("AB126").Value = 0
$ws.Range("AC126").Value = 20
$ws.Range("AD126").Value = 6
$ws.Range("AE126").Value = 10
$ws.Range("AF126").Value = 4
